$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header column
$ws.Range("H1").Value = "pos_categ_id"

# Row 2: new product
# Columns B, D and F hold numeric-looking strings (barcode / prices); a
# leading apostrophe forces Excel to keep them as text (t="s"), matching
# every other cell in this sheet, instead of silently converting them to
# numbers. The apostrophe also stamps a "quote prefix" style on the cell, so
# each such cell's style is reset back to Normal right afterwards, leaving
# the text content untouched but no stray formatting behind.
$ws.Range("B2").Value = "'7452000201667"
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Value = "Nuggets Regular Muslo Toledano 340 gr"
$ws.Range("D2").Value = "'2.79"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "EXENTO"
$ws.Range("F2").Value = "'0.00"
$ws.Range("F2").Style = "Normal"
$ws.Range("G2").Value = "https://superxtrapanama.vtexassets.com/arquivos/ids/158309-800-auto?v=637806395111300000&width=800&height=auto&aspect=true"
$ws.Range("H2").Value = "Congelados"

# Row 3: new product
$ws.Range("B3").Value = "'8711786256483"
$ws.Range("B3").Style = "Normal"
$ws.Range("C3").Value = "Papas Cong Chefs Best 2 5 Kg"
$ws.Range("D3").Value = "'5.49"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "EXENTO"
$ws.Range("F3").Value = "'0.00"
$ws.Range("F3").Style = "Normal"
$ws.Range("G3").Value = "https://superxtrapanama.vtexassets.com/arquivos/ids/191442-800-auto?v=638551093696600000&width=800&height=auto&aspect=true"
$ws.Range("H3").Value = "Congelados"

# Move the active selection to E10, matching the saved view state
$ws.Range("E10").Select()
